$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 2) of the 利润表/688560.xlsx sheet per the commit diff.

# J2 ("DATE_TYPE_CODE") is a text value that looks numeric ("002" -> "001").
# Use a leading apostrophe so it is stored as text (preserving the leading
# zero) instead of being auto-converted to the number 1.
$ws.Range("J2").Value = "'001"

$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 37211605.48
$ws.Range("P2").Value = 594618847.33
$ws.Range("Q2").Value = 545943560.27
$ws.Range("R2").Value = 50.6961659172
$ws.Range("S2").Value = 466467763.04
$ws.Range("T2").Value = 466467763.04
$ws.Range("U2").Value = 52.4046303597
$ws.Range("V2").Value = 27363435.67
$ws.Range("W2").Value = 14874021.04
$ws.Range("X2").Value = 11410689.07
$ws.Range("Y2").Value = 48126106.11
$ws.Range("Z2").Value = 47324999.04
$ws.Range("AA2").Value = 10113393.56

$ws.Range("AG2").Value = 5079285.97

$ws.Range("AP2").Value = 48.2233460403
$ws.Range("AQ2").Value = 65.804220056476
$ws.Range("AR2").Value = 45.199841221674
$ws.Range("AS2").Value = 32259905.48
$ws.Range("AT2").Value = 43.311375132797
